# "updated zsl with baltic sea"
#
# The IUCN status table gains two new columns (NA - "Not Applicable",
# inserted before NE, and RE - "Regionally Extinct", inserted before
# Total) and three new region rows (Albania, Germany, NE Atlantic),
# with the whole region list re-sorted alphabetically. Rewrite the used
# range (A1:K13) in full with the final values rather than trying to
# insert rows/columns, since several pre-existing totals (e.g. Baltic
# Sea's CR/RE split) changed in place, not just shifted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: Region, NA, NE, DD, LC, NT, VU, EN, CR, RE, Total
$header = @("Region", "NA", "NE", "DD", "LC", "NT", "VU", "EN", "CR", "RE", "Total")
for ($col = 1; $col -le $header.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $header[$col - 1]
}

# Data rows: region name followed by the 10 numeric columns above.
$data = @(
    @("Albania",            0, 0,  0,  0, 1,  0,  0, 0, 0,  1),
    @("Baltic Sea",         0, 0,  0,  0, 0,  3,  2, 2, 1,  8),
    @("Croatia",            0, 0,  8,  3, 1,  1,  5, 1, 0, 19),
    @("Europe",             0, 0, 11, 17, 7,  6,  6, 6, 0, 53),
    @("Germany",            0, 0,  0,  0, 0,  0,  0, 0, 0,  0),
    @("Ireland",            0, 7,  0,  8, 12, 4,  1, 2, 0, 34),
    @("Italy",              0, 0, 23,  7, 1,  0,  3, 4, 0, 38),
    @("Mediterranean Sea",  0, 0, 11, 10, 8, 11,  6, 16, 0, 62),
    @("NE Atlantic",        0, 0,  8, 11, 9,  2,  2, 1, 0, 33),
    @("Norway",             0, 0,  2,  0, 2,  0,  0, 1, 0,  5),
    @("Sweden",             0, 0,  0,  0, 0,  3,  1, 0, 1,  5),
    @("World",              0, 0,  2, 17, 9, 17,  9, 8, 0, 62)
)

$row = 2
foreach ($record in $data) {
    for ($col = 1; $col -le $record.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $record[$col - 1]
    }
    $row++
}
